$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '27.376.56'

$ws.Cells.Item(3, 4).Value = '1.824.55'
$ws.Cells.Item(3, 5).Value = '  +1.46%  '

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.000'
$ws.Cells.Item(4, 5).Value = '  +0.00%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '314.22'
$ws.Cells.Item(5, 5).Value = '  +1.40%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.9999'
$ws.Cells.Item(6, 5).Value = '  -0.08%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.4675'
$ws.Cells.Item(7, 5).Value = '  +4.59%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3785'
$ws.Cells.Item(8, 5).Value = '  +3.25%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.07430'
$ws.Cells.Item(9, 5).Value = '  +1.56%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.8759'
$ws.Cells.Item(10, 5).Value = '  +2.28%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '20.79'
$ws.Cells.Item(11, 5).Value = '  +0.83%  '

$ws.Cells.Item(12, 4).Value = '1.826.13'
$ws.Cells.Item(12, 5).Value = '  -4.57%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '6.689'
$ws.Cells.Item(13, 5).Value = '  +1.62%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '5.419'
$ws.Cells.Item(14, 5).Value = '  +2.64%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '93.04'
$ws.Cells.Item(15, 5).Value = '  +0.72%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.07086'
$ws.Cells.Item(16, 5).Value = '  +0.24%  '

$ws.Cells.Item(17, 5).Value = '  +0.01%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.000008792'
$ws.Cells.Item(18, 5).Value = '  +1.19%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.9998'
$ws.Cells.Item(19, 5).Value = '  -0.12%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '15.03'
$ws.Cells.Item(20, 5).Value = '  +1.37%  '

$ws.Cells.Item(21, 4).Value = '27.374.10'
$ws.Cells.Item(21, 5).Value = '  +2.13%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '5.318'
$ws.Cells.Item(22, 5).Value = '  +3.31%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '10.99'
$ws.Cells.Item(23, 5).Value = '  +2.06%  '

$ws.Cells.Item(24, 4).Value = '2.052.73'
$ws.Cells.Item(24, 5).Value = '  -3.43%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '1.941'
$ws.Cells.Item(25, 5).Value = '  -2.36%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '151.13'
$ws.Cells.Item(26, 5).Value = '  -0.46%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '2.252'
$ws.Cells.Item(27, 5).Value = '  +3.31%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '18.60'
$ws.Cells.Item(28, 5).Value = '  +1.09%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '5.339'
$ws.Cells.Item(29, 5).Value = '  +2.88%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '117.22'
$ws.Cells.Item(30, 5).Value = '  +0.48%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.08968'
$ws.Cells.Item(31, 5).Value = '  +1.88%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '0.7878'
$ws.Cells.Item(32, 5).Value = '  +5.93%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '1.194'
$ws.Cells.Item(33, 5).Value = '  +2.88%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '4.539'
$ws.Cells.Item(34, 5).Value = '  +1.91%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '2.939'
$ws.Cells.Item(35, 5).Value = '  +0.06%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.9994'
$ws.Cells.Item(36, 5).Value = '  -0.10%  '

$ws.Cells.Item(37, 5).Value = '  +1.56%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.01977'
$ws.Cells.Item(38, 5).Value = '  +0.83%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.05250'
$ws.Cells.Item(39, 5).Value = '  +1.44%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '7.293'
$ws.Cells.Item(40, 5).Value = '  +3.79%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.5355'
$ws.Cells.Item(41, 5).Value = '  +1.30%  '

$ws.Cells.Item(42, 2).Value = 'MXToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '2.897'
$ws.Cells.Item(42, 5).Value = '  +1.58%  '

$ws.Cells.Item(43, 2).Value = 'RenderToken'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '2.361'
$ws.Cells.Item(43, 5).Value = '  +20.32%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.1703'
$ws.Cells.Item(44, 5).Value = '  +1.27%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '8.655'
$ws.Cells.Item(45, 5).Value = '  +2.94%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.5093'
$ws.Cells.Item(46, 5).Value = '  -0.22%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '10.54'
$ws.Cells.Item(47, 5).Value = '  +0.35%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '105.85'
$ws.Cells.Item(48, 5).Value = '  +0.25%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.682'
$ws.Cells.Item(49, 5).Value = '  +1.27%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.9992'
$ws.Cells.Item(50, 5).Value = '  -0.11%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.06385'
$ws.Cells.Item(51, 5).Value = '  +1.34%  '
